$wb = $excel.ActiveWorkbook

function Set-CellText {
    param($ws, $r, $c, $val)
    $cell = $ws.Cells.Item($r, $c)
    if ($val -match '^-?\$?[\d,]+(\.\d+)?%?$') {
        $cell.Value2 = "'" + $val
    } else {
        $cell.Value2 = $val
    }
}

# ---- Sheet: Overall ----
$ws = $wb.Worksheets.Item("Overall")
$ws.Cells.Clear()

$header = @(
  "Share of 990 filers with government grants at risk",
  "Number of 990 filers with government grants",
  "Total government grants (`$)",
  "Size of operating surplus with government grants",
  "Size of operating surplus without government grants"
)
for ($c = 0; $c -lt $header.Length; $c++) {
    $col = $c + 1
    Set-CellText $ws 1 $col $header[$c]
}
$lastCol = $header.Length
$headerRange = $ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item(1,$lastCol))
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

$rows = @(
  @("72.14%", "718", "`$1,665,786,876", "9.73%", "-23.91%")
)
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowArr = $rows[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt $rowArr.Length; $c++) {
        $col = $c + 1
        Set-CellText $ws $rowNum $col $rowArr[$c]
    }
}

# ---- Sheet: County ----
$ws = $wb.Worksheets.Item("County")
$ws.Cells.Clear()

$header = @(
  "Geography",
  "Share of 990 filers with government grants at risk",
  "Number of 990 filers with government grants",
  "Total government grants (`$)",
  "Size of operating surplus with government grants",
  "Size of operating surplus without government grants"
)
for ($c = 0; $c -lt $header.Length; $c++) {
    $col = $c + 1
    Set-CellText $ws 1 $col $header[$c]
}
$lastCol = $header.Length
$headerRange = $ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item(1,$lastCol))
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

$rows = @(
  @("United States", "67.35%", "103,475", "`$267,700,640,005", "9.05%", "-12.83%"),
  @("Utah", "72.14%", "718", "`$1,665,786,876", "9.73%", "-23.91%"),
  @("Box Elder County", "100.00%", "6", "`$6,366,119", "-8.48%", "-83.00%"),
  @("Cache County", "74.19%", "31", "`$60,603,482", "11.68%", "-20.97%"),
  @("Carbon County", "100.00%", "8", "`$20,921,831", "8.09%", "-67.28%"),
  @("Davis County", "73.68%", "38", "`$140,229,596", "9.82%", "-16.52%"),
  @("Duchesne County", "50.00%", "2", "`$2,080,847", "5.87%", "-43.81%"),
  @("Emery County", "75.00%", "4", "`$2,616,748", "24.85%", "-37.78%"),
  @("Garfield County", "0.00%", "1", "`$101,786", "27.05%", "10.72%"),
  @("Grand County", "61.11%", "18", "`$9,590,593", "18.01%", "-10.63%"),
  @("Iron County", "100.00%", "8", "`$11,912,615", "4.08%", "-82.15%"),
  @("Juab County", "0.00%", "1", "`$3,175,040", "10.96%", "5.59%"),
  @("Kane County", "100.00%", "1", "`$170,526", "1.01%", "-21.89%"),
  @("Salt Lake County", "69.55%", "358", "`$902,347,176", "9.68%", "-22.27%"),
  @("San Juan County", "100.00%", "5", "`$30,069,443", "11.07%", "-53.04%"),
  @("Sanpete County", "33.33%", "6", "`$1,822,239", "45.57%", "26.61%"),
  @("Sevier County", "80.00%", "5", "`$3,891,632", "12.01%", "-33.79%"),
  @("Summit County", "55.00%", "40", "`$35,386,069", "14.11%", "-3.43%"),
  @("Tooele County", "75.00%", "4", "`$21,128,522", "15.11%", "-47.46%"),
  @("Uintah County", "100.00%", "2", "`$7,085,622", "8.17%", "-86.33%"),
  @("Utah County", "78.31%", "83", "`$221,809,405", "8.75%", "-41.91%"),
  @("Wasatch County", "100.00%", "6", "`$3,267,319", "8.21%", "-29.97%"),
  @("Washington County", "75.76%", "33", "`$64,070,899", "8.78%", "-53.00%"),
  @("Wayne County", "33.33%", "3", "`$3,520,279", "46.69%", "35.16%"),
  @("Weber County", "81.82%", "55", "`$113,619,088", "8.37%", "-35.78%")
)
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowArr = $rows[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt $rowArr.Length; $c++) {
        $col = $c + 1
        Set-CellText $ws $rowNum $col $rowArr[$c]
    }
}

# ---- Sheet: Congressional District ----
$ws = $wb.Worksheets.Item("Congressional District")
$ws.Cells.Clear()

$header = @(
  "Geography",
  "Share of 990 filers with government grants at risk",
  "Number of 990 filers with government grants",
  "Total government grants (`$)",
  "Size of operating surplus with government grants",
  "Size of operating surplus without government grants"
)
for ($c = 0; $c -lt $header.Length; $c++) {
    $col = $c + 1
    Set-CellText $ws 1 $col $header[$c]
}
$lastCol = $header.Length
$headerRange = $ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item(1,$lastCol))
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

$rows = @(
  @("United States", "67.35%", "103,475", "`$267,700,640,005", "9.05%", "-12.83%"),
  @("Utah", "72.14%", "718", "`$1,665,786,876", "9.73%", "-23.91%"),
  @("Congressional District 1", "70.87%", "206", "`$411,910,437", "10.03%", "-24.99%"),
  @("Congressional District 2", "73.66%", "224", "`$499,324,839", "8.56%", "-28.79%"),
  @("Congressional District 3", "73.30%", "191", "`$391,412,949", "10.42%", "-22.49%"),
  @("Congressional District 4", "69.07%", "97", "`$363,138,651", "8.75%", "-22.73%")
)
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowArr = $rows[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt $rowArr.Length; $c++) {
        $col = $c + 1
        Set-CellText $ws $rowNum $col $rowArr[$c]
    }
}

# ---- Sheet: Size ----
$ws = $wb.Worksheets.Item("Size")
$ws.Cells.Clear()

$header = @(
  "Size",
  "Share of 990 filers with government grants at risk",
  "Number of 990 filers with government grants",
  "Total government grants (`$)",
  "Size of operating surplus with government grants",
  "Size of operating surplus without government grants"
)
for ($c = 0; $c -lt $header.Length; $c++) {
    $col = $c + 1
    Set-CellText $ws 1 $col $header[$c]
}
$lastCol = $header.Length
$headerRange = $ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item(1,$lastCol))
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

$rows = @(
  @("Between `$100K and `$499K", "68.85%", "183", "`$20,324,150", "8.56%", "-18.96%"),
  @("Between `$1M and `$4.99M", "77.69%", "242", "`$372,798,560", "10.25%", "-32.08%"),
  @("Between `$500K and `$999K", "75.76%", "99", "`$29,036,290", "11.65%", "-14.11%"),
  @("Between `$5M and `$9.99M", "73.17%", "82", "`$382,453,296", "8.69%", "-74.87%"),
  @("Greater than `$10M", "62.50%", "80", "`$859,150,302", "7.56%", "-10.28%"),
  @("Less than `$100K", "59.38%", "32", "`$2,024,278", "41.39%", "-13.80%"),
  @("Total", "72.14%", "718", "`$1,665,786,876", "9.73%", "-23.91%")
)
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowArr = $rows[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt $rowArr.Length; $c++) {
        $col = $c + 1
        Set-CellText $ws $rowNum $col $rowArr[$c]
    }
}

# ---- Sheet: Subsector ----
$ws = $wb.Worksheets.Item("Subsector")
$ws.Cells.Clear()

$header = @(
  "Subsector",
  "Share of 990 filers with government grants at risk",
  "Number of 990 filers with government grants",
  "Total government grants (`$)",
  "Size of operating surplus with government grants",
  "Size of operating surplus without government grants"
)
for ($c = 0; $c -lt $header.Length; $c++) {
    $col = $c + 1
    Set-CellText $ws 1 $col $header[$c]
}
$lastCol = $header.Length
$headerRange = $ws.Range($ws.Cells.Item(1,1), $ws.Cells.Item(1,$lastCol))
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

$rows = @(
  @("Arts, Culture, and Humanities", "76.34%", "93", "`$53,901,018", "7.64%", "-19.71%"),
  @("Education (Excluding Universities)", "90.57%", "159", "`$814,755,190", "8.37%", "-85.26%"),
  @("Environment and Animals", "57.78%", "45", "`$43,451,828", "11.68%", "-8.49%"),
  @("Health (Excluding Hospitals)", "58.21%", "67", "`$186,505,861", "12.76%", "-7.62%"),
  @("Hospitals", "20.00%", "5", "`$23,597,732", "19.70%", "6.26%"),
  @("Human Services", "74.58%", "177", "`$234,857,655", "8.15%", "-26.10%"),
  @("International, Foreign Affairs", "14.29%", "7", "`$627,919", "12.14%", "2.63%"),
  @("Public, Societal Benefit", "74.36%", "39", "`$17,061,464", "9.25%", "-21.03%"),
  @("Religion Related", "25.00%", "4", "`$134,893", "23.11%", "16.60%"),
  @("Unclassified", "61.34%", "119", "`$242,105,131", "10.96%", "-14.80%"),
  @("Universities", "33.33%", "3", "`$48,788,185", "4.24%", "0.50%"),
  @("Total", "72.14%", "718", "`$1,665,786,876", "9.73%", "-23.91%")
)
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowArr = $rows[$r]
    $rowNum = $r + 2
    for ($c = 0; $c -lt $rowArr.Length; $c++) {
        $col = $c + 1
        Set-CellText $ws $rowNum $col $rowArr[$c]
    }
}

Write-Host "Edit complete"